# ============================================================================
# Applies the "ADDITIONAL SCRAPING" commit:
#  - inserts a new "Player Info" worksheet (first tab)
#  - renames MATCH_CARD_LINK -> MATCH_CODE on "ODI Batting" and "ODI Bowling"
#    and replaces the full scorecard URL values with just the numeric code
# ============================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new "Player Info" sheet in front of "ODI Batting"
# ---------------------------------------------------------------------------
$battingSheetRef = $wb.Worksheets.Item("ODI Batting")
$playerInfo = $wb.Worksheets.Add($battingSheetRef)
$playerInfo.Name = "Player Info"

# Re-fetch sheet references by name now that the sheet collection changed -
# older references can end up pointing at the wrong tab after Add().
$playerInfo = $wb.Worksheets.Item("Player Info")
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")

$piHeaders = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($i = 0; $i -lt $piHeaders.Length; $i++) {
    $playerInfo.Cells.Item(1, $i + 1).Value = $piHeaders[$i]
}

$piHeaderRange = $playerInfo.Range("A1:D1")
$piHeaderRange.Font.Bold = $true
$piHeaderRange.HorizontalAlignment = -4108
$piHeaderRange.VerticalAlignment = -4160
$piHeaderRange.Borders.LineStyle = 1

$piValues = @("4923", "Adrian Neill", "Right Handed", "Right Arm Medium")
for ($i = 0; $i -lt $piValues.Length; $i++) {
    $cell = $playerInfo.Cells.Item(2, $i + 1)
    $cell.Value = "'" + $piValues[$i]
    $cell.Style = "Normal"
}

$null = $playerInfo.Range("A1").Select()

# ---------------------------------------------------------------------------
# 2. ODI Batting: MATCH_CARD_LINK (column D) -> MATCH_CODE
# ---------------------------------------------------------------------------
$battingSheet.Range("D1").Value = "MATCH_CODE"

$battingLastRow = $battingSheet.Cells.Item($battingSheet.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $battingLastRow; $r++) {
    $cell = $battingSheet.Cells.Item($r, 4)
    $link = [string]$cell.Text
    if ($link -match "MatchCode=(\d+)") {
        $cell.Value = "'" + $matches[1]
        $cell.Style = "Normal"
    }
}

# ---------------------------------------------------------------------------
# 3. ODI Bowling: MATCH_CARD_LINK (column B) -> MATCH_CODE
# ---------------------------------------------------------------------------
$bowlingSheet.Range("B1").Value = "MATCH_CODE"

$bowlingLastRow = $bowlingSheet.Cells.Item($bowlingSheet.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $bowlingLastRow; $r++) {
    $cell = $bowlingSheet.Cells.Item($r, 2)
    $link = [string]$cell.Text
    if ($link -match "MatchCode=(\d+)") {
        $cell.Value = "'" + $matches[1]
        $cell.Style = "Normal"
    }
}

Write-Host "Sheets now:"
foreach ($s in $wb.Worksheets) {
    Write-Host " - " $s.Name
}
